$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Clear the Capacity values (column F) for the four grad dorm rows that
# didn't actually have a known capacity (rows 18-21)
$ws.Range("F18:F21").ClearContents()

# Update the active selection to F2
$ws.Range("F2").Select()
